# Updates cryptos list - refresh of price/volume columns (D, E) and the
# Polygon/OKB row swap (rows 10 and 11), per the scraped source update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "29.688.83"
$ws.Range("E2").Value = "  +3.92%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.908.18"
$ws.Range("E3").Value = "  +1.53%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.58%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.30"
$ws.Range("E5").Value = "  -0.22%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  -0.53%  "

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5165"
$ws.Range("E7").Value = "  +1.06%  "

# Row 8 - Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3973"
$ws.Range("E8").Value = "  +0.87%  "

# Row 9 - Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08496"
$ws.Range("E9").Value = "  +0.98%  "

# Row 10 - Polygon->OKB
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.59"
$ws.Range("E10").Value = "  +2.13%  "

# Row 11 - OKB->Polygon
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.121"
$ws.Range("E11").Value = "  +0.57%  "

# Row 12 - Polkadot
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.307"
$ws.Range("E12").Value = "  +0.30%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.906.69"
$ws.Range("E13").Value = "  +1.39%  "

# Row 14 - Solana
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.95"
$ws.Range("E14").Value = "  +2.13%  "

# Row 15 - Chainlink
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.341"
$ws.Range("E15").Value = "  +0.68%  "

# Row 16 - BinanceUSD
$ws.Range("E16").Value = "  -0.49%  "

# Row 17 - Litecoin
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.39"
$ws.Range("E17").Value = "  +2.12%  "

# Row 18 - ShibaInu
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001114"
$ws.Range("E18").Value = "  +0.53%  "

# Row 19 - TRON
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06740"
$ws.Range("E19").Value = "  +0.21%  "

# Row 20 - Avalanche
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.93"
$ws.Range("E20").Value = "  +1.01%  "

# Row 21 - Dai
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.64%  "

# Row 22 - Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.031"
$ws.Range("E22").Value = "  +0.92%  "

# Row 23 - WrappedBTC
$ws.Range("D23").Value = "29.708.25"
$ws.Range("E23").Value = "  +3.80%  "

# Row 24 - Cosmos
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.20"
$ws.Range("E24").Value = "  +0.51%  "

# Row 25 - Toncoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.209"
$ws.Range("E25").Value = "  -1.87%  "

# Row 26 - WrappedliquidstakedEther2.0
$ws.Range("D26").Value = "2.122.48"
$ws.Range("E26").Value = "  +1.21%  "

# Row 27 - Monero
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.17"
$ws.Range("E27").Value = "  -1.15%  "

# Row 28 - EthereumClassic
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.95"
$ws.Range("E28").Value = "  +0.82%  "

# Row 29 - LidoDAOToken
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.443"
$ws.Range("E29").Value = "  +2.47%  "

# Row 30 - BitcoinCash
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.50"
$ws.Range("E30").Value = "  +1.30%  "

# Row 31 - ImmutableX
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.075"
$ws.Range("E31").Value = "  +2.10%  "

# Row 32 - Stellar
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1055"
$ws.Range("E32").Value = "  +0.10%  "

# Row 33 - Filecoin
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.190"
$ws.Range("E33").Value = "  +6.25%  "

# Row 34 - HuobiToken
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.645"
$ws.Range("E34").Value = "  +0.93%  "

# Row 35 - VeChain
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02494"
$ws.Range("E35").Value = "  +1.16%  "

# Row 36 - Hedera
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06625"
$ws.Range("E36").Value = "  +1.22%  "

# Row 37 - FraxShare
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.100"
$ws.Range("E37").Value = "  +1.50%  "

# Row 38 - Algorand
$ws.Range("E38").Value = "  +0.43%  "

# Row 39 - InternetComputer(DFINITY)
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.215"
$ws.Range("E39").Value = "  +2.27%  "

# Row 40 - ARBITRUM
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.234"
$ws.Range("E40").Value = "  +2.81%  "

# Row 41 - TheSandbox
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6555"
$ws.Range("E41").Value = "  +1.03%  "

# Row 42 - TrustWalletToken
$ws.Range("E42").Value = "  -1.98%  "

# Row 43 - Aptos
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.32"
$ws.Range("E43").Value = "  +1.03%  "

# Row 44 - Decentraland
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6109"
$ws.Range("E44").Value = "  +0.44%  "

# Row 45 - EnergySwap
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.26"
$ws.Range("E45").Value = "  +1.32%  "

# Row 46 - PancakeSwap
$ws.Range("E46").Value = "  -0.75%  "

# Row 47 - NEARProtocol
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.063"
$ws.Range("E47").Value = "  +0.80%  "

# Row 48 - EOS
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.236"
$ws.Range("E48").Value = "  +1.26%  "

# Row 49 - Quant
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.33"
$ws.Range("E49").Value = "  +1.33%  "

# Row 50 - WEMIXTOKEN
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.156"
$ws.Range("E50").Value = "  -2.60%  "

# Row 51 - Aave
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.10"
$ws.Range("E51").Value = "  +1.05%  "
